$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing row 2.
# Each row is a hashtable keyed by column letter so we can control the
# exact order cells are populated in (affects shared-string ordering).
$data = @(
    @{ A = "user"; B = "user"; C = "Fish"; D = "Angelfish";   F = "$ 10.0"; E = "Large" },
    @{ A = "user"; B = "user"; C = "Fish"; D = "Tiger Shark"; F = "$ 12.0"; E = "Spotted" },
    @{ A = "user"; B = "user"; C = "Fish"; D = "Tiger Shark"; F = "$ 12.0"; E = "Spotless" },
    @{ A = "user"; B = "user"; C = "Fish"; D = "Goldfish";    F = "$ 12.0"; E = "Male Puppy" },
    @{ A = "user"; B = "user"; C = "Fish"; D = "Goldfish";    F = "$ 12.0"; E = "Female Puppy" },
    @{ A = "user"; B = "user"; C = "Fish"; D = "Koi";         F = "$ 12.0"; E = "Female Adult" },
    @{ A = "user"; B = "user"; C = "Fish"; D = "Koi";         F = "$ 12.0"; E = "Male Adult" }
)

$row = 3
foreach ($rowData in $data) {
    foreach ($colLetter in @("A", "B", "C", "D", "F", "E")) {
        $cell = $ws.Range($colLetter + $row)
        if ($colLetter -eq "F") {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowData[$colLetter]
    }
    $row++
}

# Adjust column D width to match the diff (stored xml width ~17.42578125).
# The COM ColumnWidth setter here snaps to whole-pixel increments, so the
# nearest achievable stored width is 17.5; use an input comfortably inside
# that pixel bucket.
$ws.Columns.Item(4).ColumnWidth = 16.7

# Update the active selection to match the diff
$ws.Range("E10").Select()
